# Apply updated crypto price/volume data to Sheet1 (cols D and E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text assignments (values that Excel will not mis-parse as numbers) ---
$ws.Range("D2").Value = "30.012.68"
$ws.Range("E2").Value = "  +9.54%  "
$ws.Range("D3").Value = "1.869.58"
$ws.Range("E3").Value = "  +6.97%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("E8").Value = "  +9.43%  "
$ws.Range("E9").Value = "  +8.39%  "
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").Value = "1.871.64"
$ws.Range("E11").Value = "  +7.14%  "
$ws.Range("E12").Value = "  +5.31%  "
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("E14").Value = "  +9.42%  "
$ws.Range("E15").Value = "  +9.74%  "
$ws.Range("E16").Value = "  +7.26%  "
$ws.Range("D17").Value = "29.981.51"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  +11.48%  "
$ws.Range("E20").Value = "  +6.05%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "2.113.50"
$ws.Range("E22").Value = "  +7.29%  "
$ws.Range("E23").Value = "  +6.51%  "
$ws.Range("E24").Value = "  +6.99%  "
$ws.Range("E25").Value = "  +7.66%  "
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("E27").Value = "  +24.06%  "
$ws.Range("E28").Value = "  +9.53%  "
$ws.Range("E29").Value = "  +5.31%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  +7.02%  "
$ws.Range("E32").Value = "  +7.76%  "
$ws.Range("E33").Value = "  +5.54%  "
$ws.Range("E34").Value = "  +7.95%  "
$ws.Range("E35").Value = "  +11.33%  "
$ws.Range("E36").Value = "  +10.24%  "
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  +16.51%  "
$ws.Range("E39").Value = "  +7.34%  "
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("E41").Value = "  +8.94%  "
$ws.Range("E42").Value = "  +6.80%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("E45").Value = "  +8.64%  "
$ws.Range("E46").Value = "  +7.96%  "
$ws.Range("E47").Value = "  +8.03%  "
$ws.Range("E48").Value = "  +4.98%  "
$ws.Range("E49").Value = "  +5.62%  "
$ws.Range("E50").Value = "  +8.33%  "
$ws.Range("E51").Value = "  +10.20%  "

# --- Price cells whose text looks like a plain number: force text format so the
#     literal string (with its exact decimal places) is preserved, then restore the
#     cell style back to the sheet default (no direct formatting) to match the source. ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4982"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2836"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06524"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07213"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6598"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.800"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007480"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.731"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.019"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.488"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.940"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.392"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.231"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08596"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.874"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6814"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.733"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9590"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.128"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4169"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.423"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.303"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3716"
$ws.Range("D51").Style = "Normal"
